$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E14:E17").Value = 3
$ws.Range("F14:F17").Value = 1
$ws.Range("G2:G5").Value = 2.629231666666667
$ws.Range("G10:G13").Value = 4.24731
$ws.Range("G14:G17").Value = 0.2859396666666667
$ws.Range("H2:H5").Value = 7.887695
$ws.Range("H10:H13").Value = 12.74193
$ws.Range("H14:H17").Value = 0.8578190000000001
$ws.Range("I2:I5").Value = 0.1414315557047068
$ws.Range("I6:I9").Value = 0.6147160060020365
$ws.Range("I10:I13").Value = 0.2284711798035388
$ws.Range("I14:I17").Value = 0.01538125848971795
$ws.Range("J2:J5").Value = 0.1414315557047067
$ws.Range("J6:J9").Value = 0.6147160060020365
$ws.Range("J10:J13").Value = 0.2284711798035388
$ws.Range("J14:J17").Value = 0.01538125848971795
$ws.Range("K2").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("K10").Value = 1
$ws.Range("K14").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.3622106666666667
$ws.Range("M3").Value = 3.642436333333333
$ws.Range("M4").Value = 2.298172333333333
$ws.Range("M5").Value = 1.433004333333334
$ws.Range("M6").Value = 0.3622106666666667
$ws.Range("M7").Value = 3.642436333333333
$ws.Range("M8").Value = 2.298172333333333
$ws.Range("M9").Value = 1.433004333333334
$ws.Range("M10").Value = 0.3622106666666667
$ws.Range("M11").Value = 3.642436333333333
$ws.Range("M12").Value = 2.298172333333333
$ws.Range("M13").Value = 1.433004333333334
$ws.Range("M14").Value = 0.3622106666666667
$ws.Range("M15").Value = 3.642436333333333
$ws.Range("M16").Value = 2.298172333333333
$ws.Range("M17").Value = 1.433004333333334
$ws.Range("N2").Value = 1.086632
$ws.Range("N4").Value = 6.894517
$ws.Range("N5").Value = 4.299013
$ws.Range("N6").Value = 1.086632
$ws.Range("N8").Value = 6.894517
$ws.Range("N9").Value = 4.299013
$ws.Range("N10").Value = 1.086632
$ws.Range("N12").Value = 6.894517
$ws.Range("N13").Value = 4.299013
$ws.Range("N14").Value = 1.086632
$ws.Range("N16").Value = 6.894517
$ws.Range("N17").Value = 4.299013
$ws.Range("O2").Value = 0.04682250814834585
$ws.Range("O3").Value = 0.4708530714096335
$ws.Range("O4").Value = 0.2970817888773835
$ws.Range("O5").Value = 0.1852426315646371
$ws.Range("O6").Value = 0.04682250814834585
$ws.Range("O7").Value = 0.4708530714096335
$ws.Range("O8").Value = 0.2970817888773835
$ws.Range("O9").Value = 0.1852426315646371
$ws.Range("O10").Value = 0.04682250814834585
$ws.Range("O11").Value = 0.4708530714096335
$ws.Range("O12").Value = 0.2970817888773835
$ws.Range("O13").Value = 0.1852426315646371
$ws.Range("O14").Value = 0.04682250814834585
$ws.Range("O15").Value = 0.4708530714096335
$ws.Range("O16").Value = 0.2970817888773835
$ws.Range("O17").Value = 0.1852426315646371
$ws.Range("P2").Value = 0.04682250814834586
$ws.Range("P3").Value = 0.4708530714096336
$ws.Range("P4").Value = 0.2970817888773835
$ws.Range("P5").Value = 0.1852426315646371
$ws.Range("P6").Value = 0.04682250814834586
$ws.Range("P7").Value = 0.4708530714096336
$ws.Range("P8").Value = 0.2970817888773835
$ws.Range("P9").Value = 0.1852426315646371
$ws.Range("P10").Value = 0.04682250814834586
$ws.Range("P11").Value = 0.4708530714096336
$ws.Range("P12").Value = 0.2970817888773835
$ws.Range("P13").Value = 0.1852426315646371
$ws.Range("P14").Value = 0.04682250814834586
$ws.Range("P15").Value = 0.4708530714096336
$ws.Range("P16").Value = 0.2970817888773835
$ws.Range("P17").Value = 0.1852426315646371
$ws.Range("Q2").Value = 0.9523357548044444
$ws.Range("Q3").Value = 9.57680895141722
$ws.Range("Q4").Value = 6.042427474257223
$ws.Range("Q5").Value = 3.767700371670556
$ws.Range("Q6").Value = 4.139217932302222
$ws.Range("Q7").Value = 41.62449970607111
$ws.Range("Q8").Value = 26.26271672559111
$ws.Range("Q9").Value = 16.37587674649778
$ws.Range("Q10").Value = 1.53842098664
$ws.Range("Q11").Value = 15.47055626293
$ws.Range("Q12").Value = 9.761050333089999
$ws.Range("Q13").Value = 6.08641363501
$ws.Range("Q14").Value = 0.1035703972897778
$ws.Range("Q15").Value = 1.041517031007889
$ws.Range("Q16").Value = 0.657138630935889
$ws.Range("Q17").Value = 0.4097527814052223
$ws.Range("R2").Value = 8.57102179324
$ws.Range("R3").Value = 86.191280562755
$ws.Range("R4").Value = 54.38184726831501
$ws.Range("R5").Value = 33.90930334503501
$ws.Range("R6").Value = 37.25296139072
$ws.Range("R7").Value = 374.62049735464
$ws.Range("R8").Value = 236.36445053032
$ws.Range("R9").Value = 147.38289071848
$ws.Range("R10").Value = 13.84578887976
$ws.Range("R11").Value = 139.23500636637
$ws.Range("R12").Value = 87.84945299781
$ws.Range("R13").Value = 54.77772271509001
$ws.Range("R14").Value = 0.9321335756080001
$ws.Range("R15").Value = 9.373653279071
$ws.Range("R16").Value = 5.914247678423001
$ws.Range("R17").Value = 3.687775032647001
$ws.Range("S2").Value = 0.006622180169416862
$ws.Range("S3").Value = 0.06659348239780384
$ws.Range("S4").Value = 0.0420167395724656
$ws.Range("S5").Value = 0.02619915356502045
$ws.Range("S6").Value = 0.02878254519994897
$ws.Range("S7").Value = 0.2894409194707216
$ws.Range("S8").Value = 0.1826209307146454
$ws.Range("S9").Value = 0.1138716106167205
$ws.Range("S10").Value = 0.01069759367801339
$ws.Range("S11").Value = 0.1075763567390789
$ws.Range("S12").Value = 0.06787462680296164
$ws.Range("S13").Value = 0.0423226025834849
$ws.Range("S14").Value = 0.0007201891009666326
$ws.Range("S15").Value = 0.007242312802029198
$ws.Range("S16").Value = 0.004569491787310851
$ws.Range("S17").Value = 0.00284926479941127
$ws.Range("T2").Value = 0.006622180169416862
$ws.Range("T3").Value = 0.06659348239780384
$ws.Range("T4").Value = 0.04201673957246559
$ws.Range("T5").Value = 0.02619915356502044
$ws.Range("T6").Value = 0.02878254519994897
$ws.Range("T7").Value = 0.2894409194707216
$ws.Range("T8").Value = 0.1826209307146454
$ws.Range("T9").Value = 0.1138716106167205
$ws.Range("T10").Value = 0.01069759367801339
$ws.Range("T11").Value = 0.1075763567390789
$ws.Range("T12").Value = 0.06787462680296164
$ws.Range("T13").Value = 0.0423226025834849
$ws.Range("T14").Value = 0.0007201891009666327
$ws.Range("T15").Value = 0.007242312802029199
$ws.Range("T16").Value = 0.004569491787310851
$ws.Range("T17").Value = 0.00284926479941127

